$d = $word.ActiveDocument
Write-Host "Paragraphs count:" $d.Paragraphs.Count
